# Insert a new record row at row 336 (weekly price-sheet refresh).
# This shifts the existing rows 336-352 down to 337-353, preserving their
# data and the date-formatted style on column D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with the new weekly record.
$ws.Range("A336").Value = 8
$ws.Range("B336").Value = "Terminal La Palmera de La Serena"
$ws.Range("C336").Value = "Coquimbo"
$ws.Range("D336").Value = 44753
$ws.Range("E336").Value = 4
$ws.Range("F336").Value = 100112032
$ws.Range("G336").Value = "Zapallo italiano"
$ws.Range("H336").Value = "Sin especificar"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 440
$ws.Range("K336").Value = 10000
$ws.Range("L336").Value = 11000
$ws.Range("M336").Value = 10500
$ws.Range("N336").Value = "`$/caja 50 unidades"
$ws.Range("O336").Value = "Región de Arica y Parinacota"
$ws.Range("P336").Value = 210
$ws.Range("Q336").Value = 50
$ws.Range("R336").Value = "Hortaliza"
